$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "fresh-water",
    "water-stone",
    "mystic-water",
    "water-gem",
    "waterium-z--held",
    "water-memory",
    "waterium-z--bag",
    "rotom-bike--water-mode",
    "scroll-of-waters",
    "water-tera-shard",
    "watercress"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $values[$i]
}

# Clear out the now-unused rows (13-35) that previously had data
$clearRange = $ws.Range("A13:A35")
$clearRange.Clear()
